$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new date value in A3, one day after A2, continuing the same date series/style
$ws.Range("A3").Value2 = $ws.Range("A2").Value2 + 1

# Copy the style/format from A2 (date format with border) down to A3
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update selection to match the target state
$ws.Range("A2:A3").Select() | Out-Null
